$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 52
Set-TextValue $ws.Range("A52") "2"
Set-TextValue $ws.Range("C52") "175_火灵鸟_Free Spirit_Rosa rugosa Thunb._20stems"
Set-TextValue $ws.Range("F52") "1.5"

# Row 53
Set-TextValue $ws.Range("C53") "177_国王日_Kings Day_Rosa rugosa Thunb._20stems"
Set-TextValue $ws.Range("F53") "14"

# Row 54
Set-TextValue $ws.Range("C54") "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"
Set-TextValue $ws.Range("F54") "10"

# Row 55
Set-TextValue $ws.Range("C55") "192_粉荔枝_Pink Ohara_Rosa rugosa Thunb._20stems"
Set-TextValue $ws.Range("F55") "7"

# Row 56
Set-TextValue $ws.Range("C56") "221_朱丽叶塔_Julieta_Rosa rugosa Thunb._10stems"
Set-TextValue $ws.Range("F56") "10"

# Row 57
Set-TextValue $ws.Range("C57") "580_腊梅黄_wax yellow_undefined_1bunch"
Set-TextValue $ws.Range("F57") "5"

# Row 58
Set-TextValue $ws.Range("C58") "582_腊梅蓝_wax blue_undefined_1bunch"
Set-TextValue $ws.Range("F58") "10"

# Row 59
Set-TextValue $ws.Range("A59") "3"
Set-TextValue $ws.Range("C59") "326_红继木_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F59") "20"

# Row 60
Set-TextValue $ws.Range("C60") "484_天鹅绒_Star of Bethlehem_undefined_1bunch"
Set-TextValue $ws.Range("F60") "10"

# Row 61
Set-TextValue $ws.Range("C61") "354_桔叶_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F61") "10"

# Row 62
Set-TextValue $ws.Range("C62") "454_蓝星花_tweedia blue_undefined_1bunch"
Set-TextValue $ws.Range("F62") "30"

# Row 63
Set-TextValue $ws.Range("C63") "505_紫罗兰紫_violet purple_undefined_1bunch"
Set-TextValue $ws.Range("F63") "10"

# Row 64
Set-TextValue $ws.Range("C64") "506_紫罗兰香槟色_violet champagne_undefined_1bunch"
Set-TextValue $ws.Range("F64") "10"

# Row 65
Set-TextValue $ws.Range("C65") "105_绣球莫奈蓝_Hydrangea Monet Blue_Hydrangea L._1stem"
Set-TextValue $ws.Range("F65") "20"

# Row 66
Set-TextValue $ws.Range("C66") "105_绣球莫奈蓝_Hydrangea Monet Blue_Hydrangea L._1stem"
Set-TextValue $ws.Range("F66") "25"

# Row 67
Set-TextValue $ws.Range("C67") "558_油画小菊_Helenium_undefined_1bunch"
Set-TextValue $ws.Range("F67") "10"

# Row 68
Set-TextValue $ws.Range("A68") "4"
Set-TextValue $ws.Range("C68") "586_洋牡丹白_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F68") "20"

# Row 69
Set-TextValue $ws.Range("C69") "589_洋牡丹香槟_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F69") "10"

# Row 70
Set-TextValue $ws.Range("C70") "651_大丽花 奶油桃子_undefined_undefined_5stems"
Set-TextValue $ws.Range("F70") "10"

# Row 71
Set-TextValue $ws.Range("C71") "507_风铃花深紫色_Canterbury Bells`ndeep purple_undefined_1bunch"

# Reset row 71 height (auto) since multi-line text may trigger custom row height
$ws.Rows.Item(71).AutoFit()

# Extend the "numbers stored as text" error-checking suppression to the newly added rows
# (mirrors the sqref="A1:L51" -> sqref="A1:L71" widening in the source workbook)
try {
    $ws.Range("A1:L71").Errors.Item(7).Ignore = $true
} catch {
}

# Update Summary!G2 aggregate string
$ws2 = $wb.Worksheets.Item("Summary")
Set-TextValue $ws2.Range("G2") "0202026271350151315142075625361010341035201420830208540445595235361550351691515258101068.5101.514107105102010103010102025102010100"
